$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$tbl = $ws.ListObjects.Item("BalanceDaily")

# New daily balance rows (bots 124-129 approx. per commit message) appended
# to the BalanceDaily table - rows 37..42 / dates 44227..44232.
$dates = 44227, 44228, 44229, 44230, 44231, 44232
$values = 0.00657683, 0.00690793, 0.00721456, 0.00738938, 0.00777784, 0.00808379

for ($i = 0; $i -lt $dates.Length; $i++) {
    $tbl.ListRows.Add() | Out-Null
    $r = 37 + $i

    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $values[$i]
    $ws.Cells.Item($r, 3).Formula = "=ROUND(IFERROR(BalanceDaily[[#This Row],[ValueLTC]]-B" + ($r - 1) + ",0),8)"
    $ws.Cells.Item($r, 4).Formula = "=BalanceDaily[[#This Row],[IncrementDaily]]/24"
}

# Match the author's final selection / scroll position on the sheet.
$ws.Range("A43").Select()

$wb.Save()
